$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) from column Q into the new column R for each data row,
# then set the R-column value for fiscal year 2021.
    $ws.Range("Q3").Copy($ws.Range("R3"))
    $ws.Range("R3").Value = 2021
    $ws.Range("Q4").Copy($ws.Range("R4"))
    $ws.Range("R4").Value = 1.7931687443515183
    $ws.Range("Q5").Copy($ws.Range("R5"))
    $ws.Range("R5").Value = 1.0977143806517458
    $ws.Range("Q6").Copy($ws.Range("R6"))
    $ws.Range("R6").Value = 2.4989281705678046
    $ws.Range("Q7").Copy($ws.Range("R7"))
    $ws.Range("R7").Value = 2.3489023398681002
    $ws.Range("Q8").Copy($ws.Range("R8"))
    $ws.Range("R8").Value = 1.8410239038543676
    $ws.Range("Q9").Copy($ws.Range("R9"))
    $ws.Range("R9").Value = 2.8382683724659588
    $ws.Range("Q10").Copy($ws.Range("R10"))
    $ws.Range("R10").Value = 1.2584206034913306
    $ws.Range("Q11").Copy($ws.Range("R11"))
    $ws.Range("R11").Value = 0.79202525610136665
    $ws.Range("Q12").Copy($ws.Range("R12"))
    $ws.Range("R12").Value = 1.7183687369364922
    $ws.Range("Q13").Copy($ws.Range("R13"))
    $ws.Range("R13").Value = 1.7860084101151579
    $ws.Range("Q14").Copy($ws.Range("R14"))
    $ws.Range("R14").Value = 1.5807090270340762
    $ws.Range("Q15").Copy($ws.Range("R15"))
    $ws.Range("R15").Value = 1.9930959157478496
    $ws.Range("Q16").Copy($ws.Range("R16"))
    $ws.Range("R16").Value = 1.0231016349164126
    $ws.Range("Q17").Copy($ws.Range("R17"))
    $ws.Range("R17").Value = 0
    $ws.Range("Q18").Copy($ws.Range("R18"))
    $ws.Range("R18").Value = 2.0091214112068791
    $ws.Range("Q19").Copy($ws.Range("R19"))
    $ws.Range("R19").Value = 2.2092990108041848
    $ws.Range("Q20").Copy($ws.Range("R20"))
    $ws.Range("R20").Value = 0.86496336159360854
    $ws.Range("Q21").Copy($ws.Range("R21"))
    $ws.Range("R21").Value = 3.5236628052020538
    $ws.Range("Q22").Copy($ws.Range("R22"))
    $ws.Range("R22").Value = 1.4678252700798498
    $ws.Range("Q23").Copy($ws.Range("R23"))
    $ws.Range("R23").Value = 0.74155920237892192
    $ws.Range("Q24").Copy($ws.Range("R24"))
    $ws.Range("R24").Value = 2.1792664589099311
    $ws.Range("Q25").Copy($ws.Range("R25"))
    $ws.Range("R25").Value = 1.5302890103825006
    $ws.Range("Q26").Copy($ws.Range("R26"))
    $ws.Range("R26").Value = 0.80351618683358383
    $ws.Range("Q27").Copy($ws.Range("R27"))
    $ws.Range("R27").Value = 2.280288974802807
    $ws.Range("Q28").Copy($ws.Range("R28"))
    $ws.Range("R28").Value = 2.3014726663297309
    $ws.Range("Q29").Copy($ws.Range("R29"))
    $ws.Range("R29").Value = 1.7358308467556451
    $ws.Range("Q30").Copy($ws.Range("R30"))
    $ws.Range("R30").Value = 2.9402079315049163
    $ws.Range("Q31").Copy($ws.Range("R31"))
    $ws.Range("R31").Value = 1.2198989923634325
    $ws.Range("Q32").Copy($ws.Range("R32"))
    $ws.Range("R32").Value = 1.1878318505232399
    $ws.Range("Q33").Copy($ws.Range("R33"))
    $ws.Range("R33").Value = 1.2537455648750642

# Update the active selection to match the target workbook state.
$ws.Range("S14").Select()
